$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new data row at row 452 (pushes old 452..481 down to 453..482) ---
$ws.Rows.Item(452).Insert()

$ws.Cells.Item(452, 1).Value = 7
$ws.Cells.Item(452, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(452, 3).Value = "Ñuble"
$ws.Cells.Item(452, 4).Value = 44610
$ws.Cells.Item(452, 5).Value = 16
$ws.Cells.Item(452, 6).Value = 100112004
$ws.Cells.Item(452, 7).Value = "Cebolla"
$ws.Cells.Item(452, 8).Value = "Morada(o)"
$ws.Cells.Item(452, 9).Value = "1a nueva(o)"
$ws.Cells.Item(452, 10).Value = 100
$ws.Cells.Item(452, 11).Value = 8000
$ws.Cells.Item(452, 12).Value = 9000
$ws.Cells.Item(452, 13).Value = 8500
$ws.Cells.Item(452, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(452, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(452, 16).Value = 472
$ws.Cells.Item(452, 17).Value = 18
$ws.Cells.Item(452, 18).Value = "Hortaliza"

# --- Insert second new data row at row 480 (pushes rows currently at 480..482 down to 481..483) ---
$ws.Rows.Item(480).Insert()

$ws.Cells.Item(480, 1).Value = 7
$ws.Cells.Item(480, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(480, 3).Value = "Ñuble"
$ws.Cells.Item(480, 4).Value = 44160
$ws.Cells.Item(480, 5).Value = 16
$ws.Cells.Item(480, 6).Value = 100112004
$ws.Cells.Item(480, 7).Value = "Cebolla"
$ws.Cells.Item(480, 8).Value = "Morada(o)"
$ws.Cells.Item(480, 9).Value = "1a (cosecha)"
$ws.Cells.Item(480, 10).Value = 75
$ws.Cells.Item(480, 11).Value = 8500
$ws.Cells.Item(480, 12).Value = 9000
$ws.Cells.Item(480, 13).Value = 8700
$ws.Cells.Item(480, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(480, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(480, 16).Value = 483
$ws.Cells.Item(480, 17).Value = 18
$ws.Cells.Item(480, 18).Value = "Hortaliza"

Write-Host "edit applied"
